$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.405.26'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.564.12'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.500'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.28%  '
$ws.Range('E9').Value = '  -1.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0591'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.787.61'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.572.10'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('E14').Value = '  -1.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.516'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.48'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.417.66'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '212.57'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.58%  '
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.25'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.20%  '
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.10'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.51'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('E24').Value = '  +2.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.76'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '14.95'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.16%  '
$ws.Range('E29').Value = '  -1.93%  '
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0469'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.20'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.373.53'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.94%  '
$ws.Range('E34').Value = '  +0.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.53'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.962'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('E37').Value = '  -0.81%  '
$ws.Range('E38').Value = '  +1.21%  '
$ws.Range('E39').Value = '  -1.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.820'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.36%  '
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.16'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.60%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.26'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.700.69'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.48'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₇0986'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.59%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0958'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.86%  '
